$d = $word.ActiveDocument

# Walk every paragraph looking for bold text ("${...}" merge terms) and
# swap Font.Bold off in favor of Font.BoldBi (w:b -> w:bCs), per the
# "remove bold font from merge terms" change.

$paraCount = $d.Paragraphs.Count

for ($p = 1; $p -le $paraCount; $p++) {
    $pRange = $d.Paragraphs.Item($p).Range
    $pStart = $pRange.Start

    $wholeBold = $pRange.Font.Bold

    if ($wholeBold -eq -1) {
        # Entire paragraph (including the paragraph mark's rPr) is bold.
        $pRange.Font.Bold = $false
        $pRange.Font.BoldBi = $true
    } else {
        # Mixed formatting inside the paragraph - walk characters and
        # flip each contiguous bold run individually so non-bold text
        # stays untouched.
        $charCount = $pRange.Characters.Count
        $spanStart = -1

        for ($i = 1; $i -le $charCount; $i++) {
            $ch = $pRange.Characters.Item($i)
            $isBold = ($ch.Font.Bold -eq -1)

            if ($isBold -and $spanStart -eq -1) {
                $spanStart = $ch.Start
            }

            $atEnd = ($i -eq $charCount)
            if ((-not $isBold) -or $atEnd) {
                if ($spanStart -ne -1) {
                    $spanEndLocal = $atEnd -and $isBold
                    if ($spanEndLocal) {
                        $localEnd = $ch.End
                    } else {
                        $localEnd = $ch.Start
                    }
                    $absStart = $pStart + $spanStart
                    $absEnd = $pStart + $localEnd
                    $boldRange = $d.Range($absStart, $absEnd)
                    $boldRange.Font.Bold = $false
                    $boldRange.Font.BoldBi = $true
                    $spanStart = -1
                }
            }
        }
    }
}
